# Insert a new "2. Visualization" slide (continuation slide, with a remark
# about pretty-printing) right after the existing "2. Visualization" slide
# (position 4) and before the "3. H-representation and V-representation"
# slide (which currently sits at position 5).
#
# Layout 2 == "Title and Content" (same layout used by the neighbouring
# slides), so the new slide gets a Title placeholder + a Content
# placeholder, just like its siblings.

$p = $ppt.ActivePresentation

$newSlide = $p.Slides.Add(5, 2)

# --- Title placeholder -------------------------------------------------
$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "2. Visualization"

# --- Content placeholder ------------------------------------------------
$body = $newSlide.Shapes.Item(2)

# Match the explicit placeholder geometry used by the neighbouring slides
# in this deck (inherited from the slide master, but pinned explicitly on
# the slide once the placeholder is edited).
$body.Left = 838200 / 12700.0
$body.Top = 1825625 / 12700.0
$body.Width = 10515600 / 12700.0
$body.Height = 4351338 / 12700.0

$tr = $body.TextFrame.TextRange
$tr.Text = "Also, extensive pretty printing and printing formats for the different classes."

# No bullet on this (single) paragraph, matching the other "intro
# sentence" paragraphs used throughout this deck.
$tr.ParagraphFormat.Bullet.Type = 0
